$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert the new "Update forward declaration header." bullet (ilvl 2) right
#    after "Forward declarations where possible. (Boost? STL?)" and give it
#    the (hidden) _GoBack bookmark that Word drops at the last edit location.
# ---------------------------------------------------------------------------
$find1 = $d.Content
$ok1 = $find1.Find.Execute("Forward declarations where possible. (Boost? STL?)", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok1) { throw "Could not find anchor paragraph for insertion" }

$srcPara = $find1.Paragraphs(1)
$srcPara.Range.InsertParagraphAfter()
$newPara = $srcPara.Next()
$newPara.Range.Text = "Update forward declaration header."

# Re-find the freshly typed text so we can stamp the _GoBack bookmark on it.
$bmRange = $d.Content
$okBm = $bmRange.Find.Execute("Update forward declaration header.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $okBm) { throw "Could not find newly inserted paragraph" }
$bmSpan = $d.Range($bmRange.Start, $bmRange.End - 1)
$d.Bookmarks.Add("_GoBack", $bmSpan)

# ---------------------------------------------------------------------------
# 2) Move <w:lastRenderedPageBreak/> off the "Performance improvements..."
#    run and onto the "Maintain a changelog." run. These runs are re-written
#    verbatim (same rsid* attributes) via InsertXML, only the rendered page
#    break marker differs.
# ---------------------------------------------------------------------------
$changelog = $d.Content
$okC = $changelog.Find.Execute("Maintain a changelog.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $okC) { throw "Could not find changelog paragraph" }
$changelog.Collapse(1)
$changelogXml = "<?xml version='1.0'?><pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body><w:p w:rsidR='00BD6916' w:rsidRPr='003342D3' w:rsidRDefault='00BD6916' w:rsidP='00BD6916'><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr><w:rPr><w:color w:val='7030A0'/></w:rPr></w:pPr><w:r w:rsidRPr='003342D3'><w:rPr><w:color w:val='7030A0'/></w:rPr><w:lastRenderedPageBreak/><w:t>Maintain a changelog.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
$changelog.InsertXML($changelogXml)

$perf = $d.Content
$okP = $perf.Find.Execute("Performance improvements in potential bottlenecks (e.g. Iterators, Scanner, PeLib,", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $okP) { throw "Could not find performance paragraph" }
$perf.Collapse(1)
$perfXml = "<?xml version='1.0'?><pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body><w:p w:rsidR='00DA2B55' w:rsidRPr='00E97FFC' w:rsidRDefault='00BD6916' w:rsidP='00E97FFC'><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr><w:rPr><w:color w:val='7030A0'/></w:rPr></w:pPr><w:r w:rsidRPr='003342D3'><w:rPr><w:color w:val='7030A0'/></w:rPr><w:t>Performance improvements in potential bottlenecks (e.g. Iterators, Scanner, PeLib,</w:t></w:r><w:r w:rsidR='006313DE' w:rsidRPr='003342D3'><w:rPr><w:color w:val='7030A0'/></w:rPr><w:t xml:space='preserve'> FindPattern,</w:t></w:r><w:r w:rsidRPr='003342D3'><w:rPr><w:color w:val='7030A0'/></w:rPr><w:t xml:space='preserve'> etc)</w:t></w:r><w:r w:rsidR='005D3274' w:rsidRPr='003342D3'><w:rPr><w:color w:val='7030A0'/></w:rPr><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
$perf.InsertXML($perfXml)

# ---------------------------------------------------------------------------
# 3) Drop the stray _GoBack bookmark that used to sit after "Add 'FreeDll'
#    API." (it moved up to the new paragraph inserted in step 1).
# ---------------------------------------------------------------------------
$freeDll = $d.Content
$okF = $freeDll.Find.Execute([string]::Format("Add {0}FreeDll{1} API.", [char]0x2018, [char]0x2019), `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $okF) { throw "Could not find FreeDll paragraph" }
$freeDll.Collapse(1)
$freeDllXml = "<?xml version='1.0'?><pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body><w:p w:rsidR='00D649A2' w:rsidRDefault='00D649A2' w:rsidP='007A1389'><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr><w:rPr><w:color w:val='7030A0'/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val='7030A0'/></w:rPr><w:t>Add " + [char]0x2018 + "FreeDll" + [char]0x2019 + " API.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
$freeDll.InsertXML($freeDllXml)

# ---------------------------------------------------------------------------
# 4) Remove the completed "Move to sub-folder to avoid header clashes etc."
#    todo item entirely (it was finished -- PeLib moved to its own folder).
# ---------------------------------------------------------------------------
$subfolder = $d.Content
$okS = $subfolder.Find.Execute("Move to sub-folder to avoid header clashes etc.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $okS) { throw "Could not find sub-folder paragraph" }
$subfolder.Paragraphs(1).Range.Delete()

Write-Output "edit complete"
